$d = $word.ActiveDocument

# Insert a brand-new, empty paragraph in front of the document's first
# paragraph (mirrors the new <w:p> the diff adds right after <w:body>).
$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphBefore()

# That new paragraph inherited the (empty) formatting context of the old
# first paragraph; the target paragraph in the diff carries no <w:pPr> at
# all, so make sure no paragraph-level alignment/etc. survived the split.
$newPara = $d.Paragraphs.First
$newPara.Format.Alignment = 0

# Give it its text. A trailing sentinel character is used so the
# collapsed "_GoBack" bookmark we add next lands cleanly inside this
# paragraph (right after "Ben Long") instead of spilling into the
# following paragraph's run.
$newPara.Range.Text = "Ben LongX"

# Re-home the "_GoBack" bookmark here, collapsed immediately after
# "Ben Long" (bookmark names are unique, so adding it here moves it off
# of its old location automatically).
$bmRange = $d.Range(8, 8)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Drop the sentinel character now that the bookmark is anchored.
$sentinel = $d.Range(8, 9)
$sentinel.Delete()
